$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Volume (F) values for rows 756-758 ---
$ws.Cells.Item(756, 6).Value = 22827.81872859
$ws.Cells.Item(757, 6).Value = 19774.57609774
$ws.Cells.Item(758, 6).Value = 7119.13938709

# --- Append new rows 760-776 ---
$newRows = @(
    @(760, 45229.45833333334, 34541.25, 34890.1, 34093.25, 34496.72, 15919.46738428),
    @(761, 45230.45833333334, 34495.65, 34728.52, 34067.72, 34660.4, 14095.02202241),
    @(762, 45231.45833333334, 34661.82, 35623.42, 34103.12, 35442.2, 24062.65560488),
    @(763, 45232.45833333334, 35442.2, 35971.45, 34328.68, 34939.68, 21258.38865655),
    @(764, 45233.45833333334, 34942.43, 34945.85, 34129.75, 34737.62, 18937.77715037),
    @(765, 45234.45833333334, 34740.63, 35277.07, 34617.56, 35091.58, 7827.22204273),
    @(766, 45235.45833333334, 35091.88, 35406.78, 34498.53, 35051.1, 9146.861512789999),
    @(767, 45236.45833333334, 35047.72, 35293.56, 34771.24, 35045.77, 10104.89174993),
    @(768, 45237.45833333334, 35054.27, 35921.03, 34534.67, 35433.57, 19161.17642952),
    @(769, 45238.45833333334, 35430.43, 36115.15, 35105.04, 35633.63, 16869.41308207),
    @(770, 45239.45833333334, 35631.88, 37980.5, 35553.71, 36698.15, 43641.33845299),
    @(771, 45240.45833333334, 36703.63, 37532.2, 36341.89, 37314.13, 22200.106037),
    @(772, 45241.45833333334, 37321.88, 37415.68, 36670.88, 37142.98, 9095.515949119999),
    @(773, 45242.45833333334, 37137.21, 37231.65, 36744.3, 37067.98, 5490.90603995),
    @(774, 45243.45833333334, 37067.6, 37432.62, 36364, 36489.44, 16552.71020945),
    @(775, 45244.45833333334, 36485.95, 36752.82, 34801.41, 35551.12, 24435.69447612),
    @(776, 45245.45833333334, 35551.63, 37965.5, 35368.04, 37881.01, 28876.22667493)
)

$dateSrc = $ws.Cells.Item(759, 1)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A carries the datetime style (same formatting as the rest of
    # the column) - copy format from the row above, then overwrite the value.
    $aCell = $ws.Cells.Item($r, 1)
    $dateSrc.Copy($aCell)
    $aCell.Value = $row[1]

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
